$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "university"

# Update A1 value to lowercase "university" (synonym search wording)
$ws.Range("A1").Value = "university"

# Move the active selection to C8, matching the author's final cursor position
$ws.Range("C8").Select()
